$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = -2
$ws.Range("F8").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("F22").Value = 4
$ws.Range("F31").Value = -1
$ws.Range("F32").Value = -4
$ws.Range("F35").Value = -5
